# Apply the SOP_MiterSaw document edits:
#  1. Main body: "Title: Prototyping Labs Supervisor" -> "Title: Prototyping Lab Supervisor"
#  2. Main body: "(if applicable)" run-merge (grammar-checker runs collapsed into a single run)
#     -- text itself stays "(if applicable)", so nothing visible changes, but make sure
#        it is present/consistent.
#  3. Header (header1.xml): "Prototyping Labs at GIX" -> "Prototyping Lab at GIX"

$d = $word.ActiveDocument

# --- 1. Fix "Prototyping Labs" -> "Prototyping Lab" in the title/supervisor line ---
$found = $d.Content.Find.Execute(
    "Title: Prototyping Labs Supervisor", $false, $false, $false, $false, $false,
    $true, 1, $false, "Title: Prototyping Lab Supervisor", 2)

# --- 2. Normalize "(if applicable)" paragraph (merges runs / clears the grammar markers) ---
$found2 = $d.Content.Find.Execute(
    "(if applicable)", $false, $false, $false, $false, $false,
    $true, 1, $false, "(if applicable)", 2)

# --- 3. Fix header text "Prototyping Labs at GIX" -> "Prototyping Lab at GIX" ---
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute(
            "Prototyping Labs at GIX", $false, $false, $false, $false, $false,
            $true, 1, $false, "Prototyping Lab at GIX", 2)
    }
}
